$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(305487936, Avihai  Kipnis: -5,-1)"
$ws.Range("B1").Value = "(313227928, Aviv  Levi: 2,6)"
$ws.Range("C1").Value = "(205807308, Sariel  Basis: 7,-1)"
$ws.Range("D1").Value = "(315891549, Raz  Halaby: -10,-4)"
$ws.Range("E1").Value = "(315060103, Dan  Mshelh: -3,4)"
$ws.Range("F1").Value = "(313925141, Elad   Amer: -1,-7)"
$ws.Range("G1").Value = "(326598423, Ron Cohen: -4,-2)"

$ws.Range("A3").Value = "cost: 431.2063753783772"
$ws.Range("A4").Value = "time: 58.02948219691102"
